$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: lunch slot moves away from 11:30 -> becomes "-"
$ws.Range("B8:F8").Value = "-"

# Row 9: time shifts 13:00 -> 12:20, content becomes "Almoço"
$ws.Range("A9").Value = "12:20"
$ws.Range("B9:F9").Value = "Almoço"

# Row 10: time shifts 13:50 -> 13:00, content stays "-"
$ws.Range("A10").Value = "13:00"

# Row 11: time shifts 14:40 -> 13:50, content stays "-"
$ws.Range("A11").Value = "13:50"

# Row 12: time shifts 15:30 -> 14:40, content becomes "-"
$ws.Range("A12").Value = "14:40"
$ws.Range("B12:F12").Value = "-"

# Row 13: time shifts 15:50 -> 15:30, content becomes "Intervalo"
$ws.Range("A13").Value = "15:30"
$ws.Range("B13:F13").Value = "Intervalo"

# Row 14 (new): 15:50, all "-"
$ws.Range("A14").Value = "15:50"
$ws.Range("B14:F14").Value = "-"

# Row 15 (new): 16:40, all "-"
$ws.Range("A15").Value = "16:40"
$ws.Range("B15:F15").Value = "-"

# Row 16 (new): 17:30, all "-"
$ws.Range("A16").Value = "17:30"
$ws.Range("B16:F16").Value = "-"

# Row 17 (new): 18:20, remaining cells empty
$ws.Range("A17").Value = "18:20"
$ws.Range("B17:F17").Value = ""
